$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 12.88963960803376
$ws.Cells.Item(2, 3).Value = 9.976687499281791
$ws.Cells.Item(2, 4).Value = 6.048845672423057
$ws.Cells.Item(2, 5).Value = 12.85754554660071
$ws.Cells.Item(2, 6).Value = 29.05216494537635
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 11).Value = 9.279679946302068
$ws.Cells.Item(2, 12).Value = 9.936603969692923
$ws.Cells.Item(2, 13).Value = 14.34080071758736
$ws.Cells.Item(2, 14).Value = 20.61222577438479
$ws.Cells.Item(2, 15).Value = 26.07061801874626
$ws.Cells.Item(3, 2).Value = 12.68735875500143
$ws.Cells.Item(3, 3).Value = 9.968956650452924
$ws.Cells.Item(3, 4).Value = 6.005592915365968
$ws.Cells.Item(3, 5).Value = 12.88524249874746
$ws.Cells.Item(3, 6).Value = 29.08328807432047
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 11).Value = 9.125311463440438
$ws.Cells.Item(3, 12).Value = 9.943880510330255
$ws.Cells.Item(3, 13).Value = 14.31423113161551
$ws.Cells.Item(3, 14).Value = 20.67265565931602
$ws.Cells.Item(3, 15).Value = 26.13503928559376
$ws.Cells.Item(4, 2).Value = 12.56411785047934
$ws.Cells.Item(4, 3).Value = 9.964333342000662
$ws.Cells.Item(4, 4).Value = 5.978386136169984
$ws.Cells.Item(4, 5).Value = 12.90391149391155
$ws.Cells.Item(4, 6).Value = 29.10890588684171
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 11).Value = 9.030839042328958
$ws.Cells.Item(4, 12).Value = 9.949666414332773
$ws.Cells.Item(4, 13).Value = 14.30005469857688
$ws.Cells.Item(4, 14).Value = 20.7115067768385
$ws.Cells.Item(4, 15).Value = 26.17937452268546
$ws.Cells.Item(5, 2).Value = 12.51420271716207
$ws.Cells.Item(5, 3).Value = 9.962480838347766
$ws.Cells.Item(5, 4).Value = 5.967139253591911
$ws.Cells.Item(5, 5).Value = 12.91193771650276
$ws.Cells.Item(5, 6).Value = 29.12098065955844
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 11).Value = 8.992468489631884
$ws.Cells.Item(5, 12).Value = 9.952356223714071
$ws.Cells.Item(5, 13).Value = 14.29481970384187
$ws.Cells.Item(5, 14).Value = 20.7277794580571
$ws.Cells.Item(5, 15).Value = 26.19864222747065
$ws.Cells.Item(6, 2).Value = 12.50593476250153
$ws.Cells.Item(6, 3).Value = 9.96217515012254
$ws.Cells.Item(6, 4).Value = 5.965262159746355
$ws.Cells.Item(6, 5).Value = 12.91329574662445
$ws.Cells.Item(6, 6).Value = 29.12308438754102
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 11).Value = 8.986106261555147
$ws.Cells.Item(6, 12).Value = 9.952822932144297
$ws.Cells.Item(6, 13).Value = 14.29398329926264
$ws.Cells.Item(6, 14).Value = 20.73050816862602
$ws.Cells.Item(6, 15).Value = 26.20191409592289
$ws.Cells.Item(7, 2).Value = 12.56344335042493
$ws.Cells.Item(7, 3).Value = 9.964308230174884
$ws.Cells.Item(7, 4).Value = 5.978235099697965
$ws.Cells.Item(7, 5).Value = 12.90401804366368
$ws.Cells.Item(7, 6).Value = 29.10906211221689
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 11).Value = 9.030320980585108
$ws.Cells.Item(7, 12).Value = 9.949701345059104
$ws.Cells.Item(7, 13).Value = 14.29998189726068
$ws.Cells.Item(7, 14).Value = 20.71172445053245
$ws.Cells.Item(7, 15).Value = 26.17962951379202
$ws.Cells.Item(8, 2).Value = 12.81973282656918
$ws.Cells.Item(8, 3).Value = 9.973996568404116
$ws.Cells.Item(8, 4).Value = 6.034068932749197
$ws.Cells.Item(8, 5).Value = 12.86675048919315
$ws.Cells.Item(8, 6).Value = 29.06154489997034
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 11).Value = 9.226418906104048
$ws.Cells.Item(8, 12).Value = 9.938839709981758
$ws.Cells.Item(8, 13).Value = 14.33119874132588
$ws.Cells.Item(8, 14).Value = 20.63270018503919
$ws.Cells.Item(8, 15).Value = 26.09183756027837
$ws.Cells.Item(9, 2).Value = 13.3270945496729
$ws.Cells.Item(9, 3).Value = 9.993954380967626
$ws.Cells.Item(9, 4).Value = 6.13826038634789
$ws.Cells.Item(9, 5).Value = 12.80685058397563
$ws.Cells.Item(9, 6).Value = 29.02003287431366
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 11).Value = 9.611278558716377
$ws.Cells.Item(9, 12).Value = 9.927971913374298
$ws.Cells.Item(9, 13).Value = 14.40916008754229
$ws.Cells.Item(9, 14).Value = 20.49153649116043
$ws.Cells.Item(9, 15).Value = 25.95766272524061
$ws.Cells.Item(10, 2).Value = 13.69902152920686
$ws.Cells.Item(10, 3).Value = 10.00917540615229
$ws.Cells.Item(10, 4).Value = 6.211389967608478
$ws.Cells.Item(10, 5).Value = 12.770859838182
$ws.Cells.Item(10, 6).Value = 29.02103937773698
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 11).Value = 9.891413959447711
$ws.Cells.Item(10, 12).Value = 9.926308670798541
$ws.Cells.Item(10, 13).Value = 14.47633927577527
$ws.Cells.Item(10, 14).Value = 20.3961547583261
$ws.Cells.Item(10, 15).Value = 25.88231120195238
$ws.Cells.Item(11, 2).Value = 13.86727499664429
$ws.Cells.Item(11, 3).Value = 10.01621531379838
$ws.Cells.Item(11, 4).Value = 6.243876626082115
$ws.Cells.Item(11, 5).Value = 12.75622353856199
$ws.Cells.Item(11, 6).Value = 29.02832593483106
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 11).Value = 10.01771776489023
$ws.Cells.Item(11, 12).Value = 9.926915129298969
$ws.Cells.Item(11, 13).Value = 14.50897808471916
$ws.Cells.Item(11, 14).Value = 20.35455442222802
$ws.Cells.Item(11, 15).Value = 25.8530876715104
$ws.Cells.Item(12, 2).Value = 13.93079220390234
$ws.Cells.Item(12, 3).Value = 10.01889727425218
$ws.Cells.Item(12, 4).Value = 6.256062959058245
$ws.Cells.Item(12, 5).Value = 12.75093043620087
$ws.Cells.Item(12, 6).Value = 29.03206479633488
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 11).Value = 10.06533809156021
$ws.Cells.Item(12, 12).Value = 9.927339842482898
$ws.Cells.Item(12, 13).Value = 14.52162974193776
$ws.Cells.Item(12, 14).Value = 20.33905744786908
$ws.Cells.Item(12, 15).Value = 25.84274903646133
$ws.Cells.Item(13, 2).Value = 13.9171223448789
$ws.Cells.Item(13, 3).Value = 10.0183189588355
$ws.Cells.Item(13, 4).Value = 6.253443612382474
$ws.Cells.Item(13, 5).Value = 12.75205931570911
$ws.Cells.Item(13, 6).Value = 29.0312160429585
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 11).Value = 10.05509216381066
$ws.Cells.Item(13, 12).Value = 9.927239713441063
$ws.Cells.Item(13, 13).Value = 14.51889210489458
$ws.Cells.Item(13, 14).Value = 20.34238362318651
$ws.Cells.Item(13, 15).Value = 25.8449432703022
$ws.Cells.Item(14, 2).Value = 13.87250488374642
$ws.Cells.Item(14, 3).Value = 10.01643564079586
$ws.Cells.Item(14, 4).Value = 6.244881540530155
$ws.Cells.Item(14, 5).Value = 12.75578307634694
$ws.Cells.Item(14, 6).Value = 29.02861391667382
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 11).Value = 10.02163994347121
$ws.Cells.Item(14, 12).Value = 9.926946167106491
$ws.Cells.Item(14, 13).Value = 14.51001313151972
$ws.Cells.Item(14, 14).Value = 20.35327434921377
$ws.Cells.Item(14, 15).Value = 25.85222251746636
$ws.Cells.Item(15, 2).Value = 13.84514798691644
$ws.Cells.Item(15, 3).Value = 10.01528413442758
$ws.Cells.Item(15, 4).Value = 6.239621860829491
$ws.Cells.Item(15, 5).Value = 12.75809645035623
$ws.Cells.Item(15, 6).Value = 29.02714752961195
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 11).Value = 10.00112106539754
$ws.Cells.Item(15, 12).Value = 9.926791734697114
$ws.Cells.Item(15, 13).Value = 14.5046123282058
$ws.Cells.Item(15, 14).Value = 20.35997856098112
$ws.Cells.Item(15, 15).Value = 25.85677605567253
$ws.Cells.Item(16, 2).Value = 13.68800155251703
$ws.Cells.Item(16, 3).Value = 10.00871761468968
$ws.Cells.Item(16, 4).Value = 6.20925088939222
$ws.Cells.Item(16, 5).Value = 12.77185126393933
$ws.Cells.Item(16, 6).Value = 29.0207004100809
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 11).Value = 9.883133008342101
$ws.Cells.Item(16, 12).Value = 9.926296386572252
$ws.Cells.Item(16, 13).Value = 14.47424750179602
$ws.Cells.Item(16, 14).Value = 20.3989093465903
$ws.Cells.Item(16, 15).Value = 25.88432281238586
$ws.Cells.Item(17, 2).Value = 13.59131230068096
$ws.Cells.Item(17, 3).Value = 10.00471853392126
$ws.Cells.Item(17, 4).Value = 6.190416943551515
$ws.Cells.Item(17, 5).Value = 12.78073384018859
$ws.Cells.Item(17, 6).Value = 29.01849285843403
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 11).Value = 9.810428399996168
$ws.Cells.Item(17, 12).Value = 9.926340972557476
$ws.Cells.Item(17, 13).Value = 14.45614728170175
$ws.Cells.Item(17, 14).Value = 20.42324959729099
$ws.Cells.Item(17, 15).Value = 25.90251705332255
$ws.Cells.Item(18, 2).Value = 13.53561446832662
$ws.Cells.Item(18, 3).Value = 10.00242930009709
$ws.Cells.Item(18, 4).Value = 6.179510907631134
$ws.Cells.Item(18, 5).Value = 12.78600628945246
$ws.Cells.Item(18, 6).Value = 29.01786606940185
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 11).Value = 9.768506757408179
$ws.Cells.Item(18, 12).Value = 9.92649497793629
$ws.Cells.Item(18, 13).Value = 14.44593259868149
$ws.Cells.Item(18, 14).Value = 20.43741794912275
$ws.Cells.Item(18, 15).Value = 25.91345758705956
$ws.Cells.Item(19, 2).Value = 13.51674351047356
$ws.Cells.Item(19, 3).Value = 10.00165609904016
$ws.Cells.Item(19, 4).Value = 6.175805836776988
$ws.Cells.Item(19, 5).Value = 12.78781952756693
$ws.Cells.Item(19, 6).Value = 29.01776433418591
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 11).Value = 9.754296389315121
$ws.Cells.Item(19, 12).Value = 9.926569195225234
$ws.Cells.Item(19, 13).Value = 14.44250796396724
$ws.Cells.Item(19, 14).Value = 20.44224408135481
$ws.Cells.Item(19, 15).Value = 25.91724353396035
$ws.Cells.Item(20, 2).Value = 13.60161426178802
$ws.Cells.Item(20, 3).Value = 10.0051431154779
$ws.Cells.Item(20, 4).Value = 6.192429455323825
$ws.Cells.Item(20, 5).Value = 12.77977136287045
$ws.Cells.Item(20, 6).Value = 29.01866132932381
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 11).Value = 9.818179016577176
$ws.Cells.Item(20, 12).Value = 9.926322948081223
$ws.Cells.Item(20, 13).Value = 14.45805383564769
$ws.Cells.Item(20, 14).Value = 20.42064110684912
$ws.Cells.Item(20, 15).Value = 25.90053100649062
$ws.Cells.Item(21, 2).Value = 13.88561593346184
$ws.Cells.Item(21, 3).Value = 10.01698838441053
$ws.Cells.Item(21, 4).Value = 6.24739959469168
$ws.Cells.Item(21, 5).Value = 12.75468255286377
$ws.Cells.Item(21, 6).Value = 29.02935166058612
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 11).Value = 10.03147167397618
$ws.Cells.Item(21, 12).Value = 9.927027102828729
$ws.Cells.Item(21, 13).Value = 14.51261323137911
$ws.Cells.Item(21, 14).Value = 20.35006853553756
$ws.Cells.Item(21, 15).Value = 25.85006467080265
$ws.Cells.Item(22, 2).Value = 14.07005385657088
$ws.Cells.Item(22, 3).Value = 10.02482353779281
$ws.Cells.Item(22, 4).Value = 6.282650053611906
$ws.Cells.Item(22, 5).Value = 12.73973879542526
$ws.Cells.Item(22, 6).Value = 29.04204665919812
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 11).Value = 10.16963768517807
$ws.Cells.Item(22, 12).Value = 9.928623776861665
$ws.Cells.Item(22, 13).Value = 14.54996982768584
$ws.Cells.Item(22, 14).Value = 20.30543792004244
$ws.Cells.Item(22, 15).Value = 25.82132378295467
$ws.Cells.Item(23, 2).Value = 13.97174303331903
$ws.Cells.Item(23, 3).Value = 10.02063338897955
$ws.Cells.Item(23, 4).Value = 6.263899169006766
$ws.Cells.Item(23, 5).Value = 12.74758169317043
$ws.Cells.Item(23, 6).Value = 29.03474976280948
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 11).Value = 10.09602320546052
$ws.Cells.Item(23, 12).Value = 9.927667947914163
$ws.Cells.Item(23, 13).Value = 14.52987881043391
$ws.Cells.Item(23, 14).Value = 20.32912191275486
$ws.Cells.Item(23, 15).Value = 25.83627496251951
$ws.Cells.Item(24, 2).Value = 13.59695708272394
$ws.Cells.Item(24, 3).Value = 10.00495113131033
$ws.Cells.Item(24, 4).Value = 6.191519841906304
$ws.Cells.Item(24, 5).Value = 12.78020598259298
$ws.Cells.Item(24, 6).Value = 29.01858316255031
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 11).Value = 9.814675341807913
$ws.Cells.Item(24, 12).Value = 9.926330697013237
$ws.Cells.Item(24, 13).Value = 14.45719128598009
$ws.Cells.Item(24, 14).Value = 20.42181986073427
$ws.Cells.Item(24, 15).Value = 25.90142740183825
$ws.Cells.Item(25, 2).Value = 13.18971984047815
$ws.Cells.Item(25, 3).Value = 9.988455321295724
$ws.Cells.Item(25, 4).Value = 6.110660836300553
$ws.Cells.Item(25, 5).Value = 12.82164571916721
$ws.Cells.Item(25, 6).Value = 29.02572806820192
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 11).Value = 9.507425693058021
$ws.Cells.Item(25, 12).Value = 9.929798799558712
$ws.Cells.Item(25, 13).Value = 14.38630708986894
$ws.Cells.Item(25, 14).Value = 20.52825598553969
$ws.Cells.Item(25, 15).Value = 25.98988623580577
